$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-04 Thursday" "2024-04-05 Friday"

Replace-Text "78×22=1716" "97×82=7954"
Replace-Text "56×26=1456" "66×98=6468"
Replace-Text "68×82=5576" "65×74=4810"
Replace-Text "60×24=1440" "68×18=1224"
Replace-Text "88×57=5016" "43×75=3225"

Replace-Text "25×47=1175" "62×18=1116"
Replace-Text "11×78=858" "64×14=896"
Replace-Text "62×52=3224" "23×60=1380"
Replace-Text "34×91=3094" "96×59=5664"
Replace-Text "76×80=6080" "71×61=4331"

Replace-Text "33×41=1353" "29×26=754"
Replace-Text "79×59=4661" "56×77=4312"
Replace-Text "91×72=6552" "18×90=1620"
Replace-Text "73×63=4599" "85×43=3655"
Replace-Text "68×57=3876" "47×26=1222"

Replace-Text "60×40=2400" "22×35=770"
Replace-Text "16×56=896" "93×27=2511"
Replace-Text "58×78=4524" "70×13=910"
Replace-Text "39×57=2223" "97×71=6887"
Replace-Text "72×55=3960" "88×65=5720"

Replace-Text "23×71=1633" "87×83=7221"
Replace-Text "86×36=3096" "96×11=1056"
Replace-Text "77×97=7469" "24×16=384"
Replace-Text "42×79=3318" "43×43=1849"
Replace-Text "11×73=803" "73×49=3577"
